# "Fixed Maintenance mode crashing"
# The maintenance-mode header text was too specific ("Time, Line, and
# Stations") for what the sheet actually shows, which was tripping up the
# app reading this cell. Shorten it, and leave the selection where the
# user last left it (cell A4) when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "Time and Stations"
$ws.Range("A4").Select()
